# Weekly refresh: insert the newest "Camote" price observation for
# "Vega Modelo de Temuco" at the top of the data block (row 87), pushing
# the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(87).Insert()

$ws.Range("A87").Value = 10
$ws.Range("B87").Value = "Vega Modelo de Temuco"
$ws.Range("C87").Value = "La Araucanía"
$ws.Range("D87").Value = "2022-08-22"
$ws.Range("E87").Value = 9
$ws.Range("F87").Value = 100114002
$ws.Range("G87").Value = "Camote"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 50
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("M87").Value = 20000
$ws.Range("N87").Value = "`$/malla 20 kilos"
$ws.Range("O87").Value = "Perú"
$ws.Range("P87").Value = 1000
$ws.Range("Q87").Value = 20
$ws.Range("R87").Value = "Hortaliza"
